$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$conv = $wb.Worksheets.Item("CONVERTION")

# Insert a new row before row 603 (shifts rows 603:641 down to 604:642,
# and the Table1 structured range auto-grows with it).
$ws1.Rows("603:603").Insert()

# Populate the newly inserted row 603.
$ws1.Range("B603").Value = "UT(0-0-18)"
$ws1.Range("D603").Value = 0.03700000000000002
$ws1.Range("G603").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# CONVERTION!F3 changes from 37 to 18 (drives G3 via the lookup table).
$conv.Range("F3").Value = 18

# Restore the active sheet/selection to match the edited workbook.
$ws1.Activate()
$ws1.Range("E606").Select()
